# Added constant distribution with no value to each user input.
#
# For every "user input" worksheet, set the Distribution Type column (F)
# to "Constant" for each parameter row, then leave the selection one row
# below the last data row (matching the post-edit cursor position), and
# finally land on the "Cost per Parameter" sheet, which becomes the
# active/selected tab when the workbook is saved.

$wb = $excel.ActiveWorkbook

# --- Incident Command (rows 2-3) ---
$ws = $wb.Worksheets.Item("Incident Command")
$ws.Range("F2").Value = "Constant"
$ws.Range("F3").Value = "Constant"
$ws.Activate()
$ws.Range("F4").Select()

# --- Characterization Sampling (rows 2-4) ---
$ws = $wb.Worksheets.Item("Characterization Sampling")
$ws.Range("F2").Value = "Constant"
$ws.Range("F3").Value = "Constant"
$ws.Range("F4").Value = "Constant"
$ws.Activate()
$ws.Range("F5").Select()

# --- Source Reduction (rows 2-4) ---
$ws = $wb.Worksheets.Item("Source Reduction")
$ws.Range("F2").Value = "Constant"
$ws.Range("F3").Value = "Constant"
$ws.Range("F4").Value = "Constant"
$ws.Activate()
$ws.Range("F5").Select()

# --- Decontamination (rows 2-6) ---
$ws = $wb.Worksheets.Item("Decontamination")
$ws.Range("F2").Value = "Constant"
$ws.Range("F3").Value = "Constant"
$ws.Range("F4").Value = "Constant"
$ws.Range("F5").Value = "Constant"
$ws.Range("F6").Value = "Constant"
# Re-apply the list validation without Delete, this drops the
# "disablePrompts" flag that Excel had set on this sheet's validation.
$dv = $ws.Range("F2:F6").Validation
$dv.Modify(3, 1, 1, "Validation_Distribution_Types")
$ws.Activate()
$ws.Range("F7").Select()

# --- Other: left untouched, no value was set for it originally ---

# --- Cost per Parameter (row 2) -- ends up the active sheet on save ---
$ws = $wb.Worksheets.Item("Cost per Parameter")
$ws.Range("F2").Value = "Constant"
$ws.Activate()
$ws.Range("F10").Select()
